$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are text strings that look numeric (e.g. "1.00", "26.966.57").
# Force them to remain text by temporarily applying a text number format, then restore
# the original cell style so no visible formatting change is introduced.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" "26.966.57"
Set-TextValue "D3" "1.560.43"
Set-TextValue "D5" "207.16"
Set-TextValue "D8" "22.14"
Set-TextValue "D12" "1.783.44"
Set-TextValue "D15" "0.519"
Set-TextValue "D16" "62.13"
Set-TextValue "D17" "26.982.10"
Set-TextValue "D18" "217.18"
Set-TextValue "D23" "9.21"
Set-TextValue "D25" "153.39"
Set-TextValue "D26" "6.61"
Set-TextValue "D27" "15.05"
Set-TextValue "D29" "1.00"
Set-TextValue "D30" "0.0469"
Set-TextValue "D31" "1.11"
Set-TextValue "D33" "1.423.34"
Set-TextValue "D36" "1.05"
Set-TextValue "D39" "0.533"
Set-TextValue "D44" "1.01"
Set-TextValue "D45" "64.77"
Set-TextValue "D47" "1.696.69"
Set-TextValue "D48" "87.39"
Set-TextValue "D50" "0.0956"

# Volume(1h) (column E) values are percentage text with surrounding spaces; assigning
# them directly keeps them as text because of the "%" sign and leading/trailing spaces.
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("E8").Value = "  +2.00%  "
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("E10").Value = "  +2.21%  "
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("E19").Value = "  +2.16%  "
$ws.Range("E20").Value = "  +2.08%  "
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("E24").Value = "  -1.32%  "
$ws.Range("E25").Value = "  -0.22%  "
$ws.Range("E26").Value = "  +0.77%  "
$ws.Range("E27").Value = "  +1.24%  "
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("E31").Value = "  +1.66%  "
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("E34").Value = "  +3.36%  "
$ws.Range("E35").Value = "  +2.80%  "
$ws.Range("E36").Value = "  +9.18%  "
$ws.Range("E37").Value = "  +1.45%  "
$ws.Range("E38").Value = "  +0.79%  "
$ws.Range("E39").Value = "  +2.65%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("E42").Value = "  +0.56%  "
$ws.Range("E43").Value = "  +2.63%  "
$ws.Range("E44").Value = "  +2.13%  "
$ws.Range("E45").Value = "  +1.75%  "
$ws.Range("E46").Value = "  +0.76%  "
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("E48").Value = "  +1.39%  "
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("E51").Value = "  -0.11%  "
